$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.667.38'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.16%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.290.53'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.00%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.25%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '112.83'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +16.80%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '268.50'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.45%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.23%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.29%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.614'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +1.42%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '47.32'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +4.45%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0936'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.15%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.51'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +8.84%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.107'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.41%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.55'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +2.86%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.630.75'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.10%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.845'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.45%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.286.44'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.40%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.560.00'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.03%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0000110'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.63%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.55'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +6.15%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.31'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.52%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.51'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +4.10%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '232.74'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.49'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +3.34%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.83'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +13.80%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.39'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +2.53%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '42.71'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +5.56%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.30%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.38%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '176.31'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.62%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '21.64'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.15%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0922'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +4.32%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.49'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +2.25%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.79%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.70'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +7.60%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.93%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.02%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.81'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +12.94%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.42'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +4.07%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '73.36'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +15.06%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +2.50%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '13.47'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +10.25%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.43'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.14%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.95'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +13.70%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.75'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.51%  '
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '102.12'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +3.98%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0998'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -2.14%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +3.15%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +4.27%  '

Write-Host "Applied all cryptos updates"